$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "09768edd95a8b219f10218dc50a94417"
$ws.Range("B17").Value = "1ccfc1ec97dfed9f35c1ed5011b1cea9"
$ws.Range("B34").Value = "62d0f019011e1e35afb4da08a53861dd"
$ws.Range("B126").Value = "f8e2c2e76e50c47fd884009976743833"
$ws.Range("B136").Value = "eb7b0979e989c558249db2170fe6a48d"
$ws.Range("B159").Value = "4749c882ce4f82f5ec89fee91ecc415c"
$ws.Range("B162").Value = "d9cbdf45e33118bc240620a3976be092"
$ws.Range("B169").Value = "4da83de0fa8baa0c3e34ef948fa497bf"
$ws.Range("B180").Value = "ae42a0af0e2092a422639ad4d71db265"
$ws.Range("B183").Value = "477b146f8b21754abe9e6418d07f97ae"
$ws.Range("B200").Value = "875decfdb4d3f6746c65a89f45459306"
$ws.Range("B213").Value = "618db607106c4c865cbafcf8156b579a"
$ws.Range("B228").Value = "5b813c348de89f8832b3df7554abeb70"
$ws.Range("B276").Value = "8c6e2b75376b8490b816902250befb49"
$ws.Range("B284").Value = "afc91a4d0896544a39504d970bebe301"
$ws.Range("B305").Value = "41f7a08e5604f7733de62b092e819c2d"
$ws.Range("B342").Value = "1eb832b6afed5fa4baf694d891211e50"
$ws.Range("B467").Value = "3c75af0a389448ba653dbb96b057f85d"
$ws.Range("B468").Value = "e1e4b714dddf2e3deb6075c4d94ffcf9"
$ws.Range("B509").Value = "4670f7f253d8abe8a660119fd708e885"
$ws.Range("B510").Value = "bcf10a301975099317a3671d48f56727"
$ws.Range("B516").Value = "0f2b68cdf56bae47118f70f03e78d2f5"
$ws.Range("B527").Value = "47b1b203b6ab8a70b7b10583d0108c5b"
$ws.Range("B562").Value = "500fec36363758d7e706ee1f3a320cbd"
$ws.Range("B584").Value = "90e9978e5fac4cdc1c413f6cc4049a3c"
$ws.Range("B628").Value = "a619418188285d32ee4afa2a1af3c1ad"
$ws.Range("B639").Value = "eff5797203762a41ac372a1640233c11"
$ws.Range("B692").Value = "87f7d8c8d5f14748512c9245c79f6ea6"
$ws.Range("B697").Value = "e992428de39ad6cc52cb72f089587295"
$ws.Range("B712").Value = "c73244e4d02da93b2f5418460dd36c9d"
$ws.Range("B715").Value = "d174fa8fbca0c777f41402c2571309ad"
$ws.Range("B727").Value = "c5ee5882e46f01af84add9b219ddf0c2"
$ws.Range("B741").Value = "1f9b18a75e7137204200fd2e581624f2"
$ws.Range("B831").Value = "3cad1c31d6cda35f1ce8b17cbb9cfdb9"
$ws.Range("B842").Value = "e08d817cc6a46610a3b5f893585aa94e"
$ws.Range("B847").Value = "b102e7c044aa28ec0c96f4f071d794ab"
$ws.Range("B874").Value = "376b400271a9aac22e19182e385681ae"
